$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189 - Al Quwarah / القوارة (منطقة القصيم / وسط المملكة)
$ws.Cells.Item(189, 1).Value = "Al Quwarah"
$ws.Cells.Item(189, 2).Value = "Al Quwarah"
$ws.Cells.Item(189, 3).Value = "القوارة"
$ws.Cells.Item(189, 4).Value = 26.770962999999998
$ws.Cells.Item(189, 5).Value = 43.473886999999998
$ws.Cells.Item(189, 6).Value = "منطقة القصيم"
$ws.Cells.Item(189, 7).Value = "وسط المملكة"

# Row 190 - Meegowa / ميقوع (منطقة الجوف / شمال المملكة)
$ws.Cells.Item(190, 1).Value = "Meegowa"
$ws.Cells.Item(190, 2).Value = "Meegowa"
$ws.Cells.Item(190, 3).Value = "ميقوع"
$ws.Cells.Item(190, 4).Value = 29.814852999999999
$ws.Cells.Item(190, 5).Value = 38.918719000000003
$ws.Cells.Item(190, 6).Value = "منطقة الجوف"
$ws.Cells.Item(190, 7).Value = "شمال المملكة"

# Match the formatting (thin border on all sides, same as the rest of the table)
# already used on the preceding data row by copying its format onto the new rows.
$ws.Range("A188:G188").Copy()
$ws.Range("A189:G190").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Refresh the selection / frozen-pane scroll position to cover the new extent,
# mirroring the view state Excel persists after editing near the end of the list.
[void]$ws.Range("A179").Select()
$excel.ActiveWindow.ScrollRow = 179
[void]$ws.Range("A1:G190").Select()
